$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25-58 down to 26-59),
# mirroring Excel's own "insert row, pushing cells down" behaviour so that
# formatting (e.g. the date style on column D) is inherited correctly.
$ws.Rows(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value = 7
$ws.Cells.Item(25, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(25, 3).Value = "Ñuble"
$ws.Cells.Item(25, 4).Value = 45210
$ws.Cells.Item(25, 5).Value = 16
$ws.Cells.Item(25, 6).Value = 300000000
$ws.Cells.Item(25, 7).Value = "Espárragos"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 1200
$ws.Cells.Item(25, 12).Value = 1300
$ws.Cells.Item(25, 13).Value = 1250
$ws.Cells.Item(25, 14).Value = "`$/kilo"
$ws.Cells.Item(25, 15).Value = "Región de Ñuble"
$ws.Cells.Item(25, 16).Value = 1250
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"
